$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price-column cells whose new values would
# otherwise be auto-converted to numbers by Excel (losing exact formatting,
# e.g. "38.00" -> 38). All data cells in this sheet are plain text (inlineStr).
$ws.Range("D2").Value = "57.026.93"
$ws.Range("E2").Value = "  +6.55%  "
$ws.Range("D3").Value = "3.236.18"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.48"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.63"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  +4.54%  "
$ws.Range("D8").Value = "3.233.31"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.615"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.12"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0959"
$ws.Range("E12").Value = "  +9.67%  "
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "3.743.77"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.19"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.10"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "3.223.20"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.97"
$ws.Range("E19").Value = "  +3.50%  "
$ws.Range("D20").Value = "56.816.26"
$ws.Range("E20").Value = "  +6.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.35"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("E22").Value = "  +7.51%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "297.22"
$ws.Range("E24").Value = "  +9.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.09"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.16"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.87"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  -5.02%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.42"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "38.00"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0484"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.12"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.73"
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").Value = "  +7.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "134.32"
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.99"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.06"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.282"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.17"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.156.93"
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.03"
$ws.Range("E50").Value = "  +23.16%  "
$ws.Range("E51").Value = "  -2.54%  "
